# lsh_coding.xlsx edits:
#  - lsh_unit_categories: correct the unit category for "Fv-G2 BM Göngu" (row 23,
#    B/C columns) and append a new row for "Fv-G3 BM Göngu" (row 26)
#  - lsh_sheet_names: append the new sheet-name entry "Takmörkun meðferðar" (row 14)
#  - switch the active sheet/tab from clinical_assessment_categories to
#    lsh_sheet_names, and update the selections on both affected sheets

$wb = $excel.ActiveWorkbook

# --- lsh_unit_categories (sheet3) ---------------------------------------
$wsUnits = $wb.Worksheets.Item("lsh_unit_categories")
$wsUnits.Activate()

# Row 23 ("Fv-G2 BM Göngu") was miscoded as Göngudeild/outpatient_clinic;
# correct it to Bráðamóttaka/emergency_room.
$wsUnits.Range("B23").Value = "Bráðamóttaka"
$wsUnits.Range("C23").Value = "emergency_room"

# New row 26 for "Fv-G3 BM Göngu", coded the same way as row 23 now is.
$wsUnits.Range("A26").Value = "Fv-G3 BM Göngu"
$wsUnits.Range("B26").Value = "Bráðamóttaka"
$wsUnits.Range("C26").Value = "emergency_room"
$wsUnits.Range("D26").Value = "home"
$wsUnits.Range("E26").Value = 1

# Scroll the view down a bit and leave the new row selected.
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$wsUnits.Range("B26").Select()

# --- lsh_sheet_names (sheet6) -------------------------------------------
$wsNames = $wb.Worksheets.Item("lsh_sheet_names")
$wsNames.Activate()

# New row 14 with the most recently added sheet name.
$wsNames.Range("A14").Value = "Takmörkun meðferðar"

$wsNames.Range("D18").Select()
